$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 values to "custom accuracy" (2 decimal places) ---
$ws.Range("C5").Value = 13.91
$ws.Range("D5").Value = 1.15
$ws.Range("E5").Value = 40.75
$ws.Range("F5").Value = 33.15
$ws.Range("G5").Value = 14.36
$ws.Range("H5").Value = 55.56
$ws.Range("J5").Value = 10.17
$ws.Range("K5").Value = 14.69
$ws.Range("L5").Value = 16.48
$ws.Range("M5").Value = 17.28
$ws.Range("O5").Value = 14.66
$ws.Range("P5").Value = 20.77
$ws.Range("Q5").Value = 12.45
$ws.Range("R5").Value = 0.69
$ws.Range("S5").Value = 0.6
$ws.Range("T5").Value = 215.86
$ws.Range("V5").Value = 13.53
$ws.Range("W5").Value = 27.45
$ws.Range("X5").Value = 14.49
$ws.Range("Z5").Value = 27.76
$ws.Range("AA5").Value = 11.95
$ws.Range("AB5").Value = 10.63
$ws.Range("AC5").Value = 12.49
$ws.Range("AD5").Value = 17.19
$ws.Range("AF5").Value = 50.6
$ws.Range("AG5").Value = 7.59
$ws.Range("AH5").Value = 16.92

# --- Remove row 6 (dimension shrinks to A1:AH5) ---
$ws.Rows.Item(6).Delete()

# --- Narrow a subset of data columns from width 8 to width 7 ---
# (ColumnWidth setter is expressed in Excel "characters"; 6.14 round-trips
# to a stored raw width of 7 for this sheet's default font/MDW.)
$ws.Columns.Item(3).ColumnWidth = 6.14    # C
$ws.Columns.Item(7).ColumnWidth = 6.14    # G
$ws.Columns.Item(10).ColumnWidth = 6.14   # J
$ws.Columns.Item(11).ColumnWidth = 6.14   # K
$ws.Columns.Item(15).ColumnWidth = 6.14   # O
$ws.Columns.Item(17).ColumnWidth = 6.14   # Q
$ws.Columns.Item(22).ColumnWidth = 6.14   # V
$ws.Columns.Item(24).ColumnWidth = 6.14   # X
$ws.Columns.Item(27).ColumnWidth = 6.14   # AA
$ws.Columns.Item(28).ColumnWidth = 6.14   # AB
$ws.Columns.Item(29).ColumnWidth = 6.14   # AC
